$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for Perejil at
# "Terminal Hortofrutícola Agro Chillán". Insert a new row at position 21,
# which shifts the existing rows 21-46 down to 22-47 (and grows the used
# range from A1:R46 to A1:R47), then populate the new row with its data.
$ws.Rows("21:21").Insert()

$ws.Cells.Item(21, 1).Value = 7
$ws.Cells.Item(21, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(21, 3).Value = "Ñuble"
$ws.Cells.Item(21, 4).Value = 45036
$ws.Cells.Item(21, 5).Value = 16
$ws.Cells.Item(21, 6).Value = 100112044
$ws.Cells.Item(21, 7).Value = "Perejil"
$ws.Cells.Item(21, 8).Value = "Sin especificar"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 250
$ws.Cells.Item(21, 11).Value = 1500
$ws.Cells.Item(21, 12).Value = 1500
$ws.Cells.Item(21, 13).Value = 1500
$ws.Cells.Item(21, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(21, 15).Value = "Región del Maule"
$ws.Cells.Item(21, 16).Value = 1500
$ws.Cells.Item(21, 17).Value = 1
$ws.Cells.Item(21, 18).Value = "Hortaliza"
